# Insert a new data row at row 172 (pushing existing rows 172:285 down to 173:286)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(172).Insert()

$ws.Range("A172").Value = 6
$ws.Range("B172").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C172").Value = "Metropolitana"
$ws.Range("D172").Value = 44518
$ws.Range("E172").Value = 13
$ws.Range("F172").Value = 100112030
$ws.Range("G172").Value = "Poroto granado"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 400
$ws.Range("K172").Value = 33000
$ws.Range("L172").Value = 35000
$ws.Range("M172").Value = 34150
$ws.Range("N172").Value = "$/malla 25 kilos"
$ws.Range("O172").Value = "Perú"
$ws.Range("P172").Value = 1366
$ws.Range("Q172").Value = 25
$ws.Range("R172").Value = "Hortaliza"
